$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 47
$ws1.Range("F4").Value = 615
$ws1.Range("F5").Value = 140
$ws1.Range("F6").Value = 9275
$ws1.Range("F7").Value = 836
$ws1.Range("F9").Value = 1186
$ws1.Range("F10").Value = 1082
$ws1.Range("F12").Value = 71
$ws1.Range("F15").Value = 384
$ws1.Range("F17").Value = 247
$ws1.Range("F18").Value = 1220

# Sheet "全部类型" (All types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 47
$ws4.Range("F6").Value = 615
$ws4.Range("F7").Value = 140
$ws4.Range("F8").Value = 9275
$ws4.Range("F9").Value = 836
$ws4.Range("F11").Value = 1186
$ws4.Range("F12").Value = 1082
$ws4.Range("F14").Value = 71
$ws4.Range("F17").Value = 384
$ws4.Range("F19").Value = 247
$ws4.Range("F20").Value = 1220
